$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("D20").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("C25").Value = -0.007000000000000006
$ws.Range("E25").Value = -0.02399999999999991
$ws.Range("C28").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("E31").Value = -0.01400000000000001
$ws.Range("F31").Value = -0.04000000000000004
$ws.Range("C32").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("H38").Value = 0.03900000000000003
$ws.Range("L38").Value = 0.01800000000000002
$ws.Range("F43").Value = -0.01200000000000001
$ws.Range("G43").Value = 0.02900000000000003
$ws.Range("G44").Value = 0.03200000000000003
$ws.Range("H45").Value = 0
$ws.Range("C46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("F48").Value = 0.02900000000000003
$ws.Range("G48").Value = 0.03300000000000003
$ws.Range("E49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("G54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("E56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("G59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("E60").Value = -0.02999999999999992
$ws.Range("L61").Value = 0
$ws.Range("G62").Value = 0.08699999999999997
$ws.Range("D63").Value = -0.02100000000000002
$ws.Range("G63").Value = 0.01700000000000002
$ws.Range("J65").Value = 0
$ws.Range("C66").Value = -0.01800000000000002
$ws.Range("J68").Value = 0.149
$ws.Range("H69").Value = 0.07300000000000001
$ws.Range("L70").Value = -0.02200000000000002
$ws.Range("E71").Value = -0.01900000000000002
$ws.Range("I71").Value = 0.02899999999999991
$ws.Range("L73").Value = -0.03400000000000003
$ws.Range("L75").Value = -0.01400000000000001
$ws.Range("C76").Value = 0.03600000000000003
$ws.Range("J76").Value = 0.03700000000000003
$ws.Range("E77").Value = -0.06999999999999995
$ws.Range("B78").Value = 0
$ws.Range("H79").Value = 0.07100000000000001
$ws.Range("I79").Value = 0.03500000000000003
$ws.Range("E80").Value = -0.02500000000000002
$ws.Range("G81").Value = 0.129
$ws.Range("H81").Value = 0.02400000000000002
$ws.Range("G82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("C83").Value = -0.04900000000000004
$ws.Range("F84").Value = 0.145
$ws.Range("H84").Value = 0.168
$ws.Range("K85").Value = -0.02200000000000002
$ws.Range("E86").Value = -0.05199999999999994
$ws.Range("I86").Value = -0.001000000000000001
$ws.Range("D87").Value = -0.03900000000000003
$ws.Range("E87").Value = -0.02999999999999992
$ws.Range("B88").Value = 0.118
$ws.Range("J88").Value = 0.1919999999999999
$ws.Range("C89").Value = 0.01300000000000001
$ws.Range("J89").Value = 0.08299999999999996
$ws.Range("H90").Value = 0.065
$ws.Range("F91").Value = 0.02799999999999991
$ws.Range("G91").Value = 0.06399999999999995
$ws.Range("C93").Value = -0.003000000000000003
$ws.Range("G93").Value = 0.06000000000000005
$ws.Range("G95").Value = 0.132
$ws.Range("G96").Value = 0.02800000000000002
$ws.Range("H97").Value = 0.155
$ws.Range("D98").Value = 0.03999999999999998
$ws.Range("H98").Value = 0.05199999999999999
$ws.Range("J99").Value = 0.104
$ws.Range("J100").Value = 0.142
